$wb = $excel.ActiveWorkbook

# Duplicate an existing sheet so the new tab inherits the same background
# picture / zoom settings as the rest of the workbook, then move it to the
# end and rename + repopulate it as "LAST".
$src = $wb.Worksheets.Item("TEXT")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "LAST"
$ws.Cells.Clear()

$ws.Range("A1:D4").NumberFormat = "@"

$ws.Range("A1").Value = "NAME"
$ws.Range("A2").Value = "TRANS"

$ws.Range("B1").Value = "PA"
$ws.Range("B2").Value = "25"

$ws.Range("C1").Value = "TA"
$ws.Range("C2").Value = "120"

$ws.Range("D1").Value = "GAMMED"
$ws.Range("D2").Value = "1000"

$ws.Range("A3").Value = "NAME"
$ws.Range("B3").Value = "PA"
$ws.Range("C3").Value = "TA"
$ws.Range("D3").Value = "GAMMED"

$ws.Range("A4").Value = "DISTR"
$ws.Range("B4").Value = "16"
$ws.Range("C4").Value = "90"
$ws.Range("D4").Value = "1000"

$ws.Range("D4").Select() | Out-Null
